# Clear out the three existing data rows (BR50050, BR50051, BR50052 — the
# last three records in the sheet) so the table body itself is emptied,
# then extend the (now-blank) table down with more empty rows, matching the
# date-formatted style that column L (Date Added) already carries.
#
# Net effect matches the upstream commit: rows 11-13 lose their content,
# rows 14-52 appear as blank placeholder rows (all carrying the date
# number-format on column L, inherited via copy/paste-format from the
# existing L11 cell), and the now-unused shared strings for those three
# removed records drop out of the shared string table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the contents of the last 3 populated rows.
[void]$ws.Range("A11:AO13").ClearContents()

# Carry the date number-format (already applied to L11:L13) down through
# L52 so the newly-extended rows pick up the same style index as the
# cleared rows above them.
[void]$ws.Range("L11").Copy()
[void]$ws.Range("L14:L52").PasteSpecial(-4122)

# Leave the selection where the edit session ended.
[void]$ws.Range("D26").Select()
